$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Cells.Item(2, 4)
$c.NumberFormat = "@"
$c.Value = '62.572.09'
$c.Style = "Normal"
$ws.Cells.Item(2, 5).Value = '  -1.32%  '

$c = $ws.Cells.Item(3, 4)
$c.NumberFormat = "@"
$c.Value = '3.016.12'
$c.Style = "Normal"
$ws.Cells.Item(3, 5).Value = '  -1.52%  '

$c = $ws.Cells.Item(4, 4)
$c.NumberFormat = "@"
$c.Value = '1.00'
$c.Style = "Normal"
$ws.Cells.Item(4, 5).Value = '  -0.03%  '

$c = $ws.Cells.Item(5, 4)
$c.NumberFormat = "@"
$c.Value = '584.96'
$c.Style = "Normal"
$ws.Cells.Item(5, 5).Value = '  -1.26%  '

$c = $ws.Cells.Item(6, 4)
$c.NumberFormat = "@"
$c.Value = '146.49'
$c.Style = "Normal"
$ws.Cells.Item(6, 5).Value = '  -5.06%  '

$ws.Cells.Item(7, 5).Value = '  +0.00%  '

$c = $ws.Cells.Item(8, 4)
$c.NumberFormat = "@"
$c.Value = '0.528'
$c.Style = "Normal"
$ws.Cells.Item(8, 5).Value = '  -2.11%  '

$c = $ws.Cells.Item(9, 4)
$c.NumberFormat = "@"
$c.Value = '3.012.08'
$c.Style = "Normal"
$ws.Cells.Item(9, 5).Value = '  -1.64%  '

$ws.Cells.Item(10, 5).Value = '  -4.07%  '

$c = $ws.Cells.Item(11, 4)
$c.NumberFormat = "@"
$c.Value = '5.80'
$c.Style = "Normal"
$ws.Cells.Item(11, 5).Value = '  -0.52%  '

$c = $ws.Cells.Item(12, 4)
$c.NumberFormat = "@"
$c.Value = '0.464'
$c.Style = "Normal"
$ws.Cells.Item(12, 5).Value = '  +2.86%  '

$c = $ws.Cells.Item(13, 4)
$c.NumberFormat = "@"
$c.Value = '0.0000230'
$c.Style = "Normal"
$ws.Cells.Item(13, 5).Value = '  -3.07%  '

$c = $ws.Cells.Item(14, 4)
$c.NumberFormat = "@"
$c.Value = '34.63'
$c.Style = "Normal"
$ws.Cells.Item(14, 5).Value = '  -6.08%  '

$ws.Cells.Item(15, 5).Value = '  +2.30%  '

$c = $ws.Cells.Item(16, 4)
$c.NumberFormat = "@"
$c.Value = '3.510.20'
$c.Style = "Normal"
$ws.Cells.Item(16, 5).Value = '  -1.64%  '

$c = $ws.Cells.Item(17, 4)
$c.NumberFormat = "@"
$c.Value = '7.14'
$c.Style = "Normal"
$ws.Cells.Item(17, 5).Value = '  -0.34%  '

$c = $ws.Cells.Item(18, 4)
$c.NumberFormat = "@"
$c.Value = '62.520.65'
$c.Style = "Normal"
$ws.Cells.Item(18, 5).Value = '  -1.36%  '

$c = $ws.Cells.Item(19, 4)
$c.NumberFormat = "@"
$c.Value = '3.014.57'
$c.Style = "Normal"
$ws.Cells.Item(19, 5).Value = '  -1.69%  '

$c = $ws.Cells.Item(20, 4)
$c.NumberFormat = "@"
$c.Value = '460.12'
$c.Style = "Normal"
$ws.Cells.Item(20, 5).Value = '  -6.15%  '

$c = $ws.Cells.Item(21, 4)
$c.NumberFormat = "@"
$c.Value = '14.04'
$c.Style = "Normal"
$ws.Cells.Item(21, 5).Value = '  -2.46%  '

$c = $ws.Cells.Item(22, 4)
$c.NumberFormat = "@"
$c.Value = '0.692'
$c.Style = "Normal"
$ws.Cells.Item(22, 5).Value = '  -1.96%  '

$c = $ws.Cells.Item(23, 4)
$c.NumberFormat = "@"
$c.Value = '7.45'
$c.Style = "Normal"
$ws.Cells.Item(23, 5).Value = '  -1.26%  '

$c = $ws.Cells.Item(24, 4)
$c.NumberFormat = "@"
$c.Value = '81.74'
$c.Style = "Normal"
$ws.Cells.Item(24, 5).Value = '  -0.25%  '

$ws.Cells.Item(25, 2).Value = 'Fetch.AI'
$ws.Cells.Item(25, 3).Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$c = $ws.Cells.Item(25, 4)
$c.NumberFormat = "@"
$c.Value = '2.23'
$c.Style = "Normal"
$ws.Cells.Item(25, 5).Value = '  -8.50%  '

$ws.Cells.Item(26, 2).Value = 'InternetComputer(DFINITY)'
$ws.Cells.Item(26, 3).Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$c = $ws.Cells.Item(26, 4)
$c.NumberFormat = "@"
$c.Value = '12.37'
$c.Style = "Normal"
$ws.Cells.Item(26, 5).Value = '  -3.80%  '

$ws.Cells.Item(27, 2).Value = 'RenderToken'
$ws.Cells.Item(27, 3).Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$c = $ws.Cells.Item(27, 4)
$c.NumberFormat = "@"
$c.Value = '10.08'
$c.Style = "Normal"
$ws.Cells.Item(27, 5).Value = '  -6.11%  '

$ws.Cells.Item(28, 2).Value = 'Dai'
$ws.Cells.Item(28, 3).Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$c = $ws.Cells.Item(28, 4)
$c.NumberFormat = "@"
$c.Value = '1.00'
$c.Style = "Normal"
$ws.Cells.Item(28, 5).Value = '  +0.15%  '

$ws.Cells.Item(29, 5).Value = '  -0.15%  '

$c = $ws.Cells.Item(30, 4)
$c.NumberFormat = "@"
$c.Value = '2.62'
$c.Style = "Normal"
$ws.Cells.Item(30, 5).Value = '  -2.45%  '

$c = $ws.Cells.Item(31, 4)
$c.NumberFormat = "@"
$c.Value = '7.07'
$c.Style = "Normal"
$ws.Cells.Item(31, 5).Value = '  -4.23%  '

$c = $ws.Cells.Item(32, 4)
$c.NumberFormat = "@"
$c.Value = '2.10'
$c.Style = "Normal"
$ws.Cells.Item(32, 5).Value = '  -5.23%  '

$c = $ws.Cells.Item(33, 4)
$c.NumberFormat = "@"
$c.Value = '28.10'
$c.Style = "Normal"
$ws.Cells.Item(33, 5).Value = '  +2.84%  '

$c = $ws.Cells.Item(34, 4)
$c.NumberFormat = "@"
$c.Value = '0.110'
$c.Style = "Normal"
$ws.Cells.Item(34, 5).Value = '  -1.65%  '

$c = $ws.Cells.Item(35, 4)
$c.NumberFormat = "@"
$c.Value = '0.0₃0808'
$c.Style = "Normal"
$ws.Cells.Item(35, 5).Value = '  -1.76%  '

$ws.Cells.Item(36, 5).Value = '  -3.17%  '

$c = $ws.Cells.Item(37, 4)
$c.NumberFormat = "@"
$c.Value = '5.78'
$c.Style = "Normal"
$ws.Cells.Item(37, 5).Value = '  -3.33%  '

$c = $ws.Cells.Item(38, 4)
$c.NumberFormat = "@"
$c.Value = '2.13'
$c.Style = "Normal"
$ws.Cells.Item(38, 5).Value = '  -4.49%  '

$c = $ws.Cells.Item(39, 4)
$c.NumberFormat = "@"
$c.Value = '50.45'
$c.Style = "Normal"
$ws.Cells.Item(39, 5).Value = '  -0.32%  '

$c = $ws.Cells.Item(40, 4)
$c.NumberFormat = "@"
$c.Value = '9.17'
$c.Style = "Normal"
$ws.Cells.Item(40, 5).Value = '  -0.80%  '

$c = $ws.Cells.Item(41, 4)
$c.NumberFormat = "@"
$c.Value = '2.92'
$c.Style = "Normal"
$ws.Cells.Item(41, 5).Value = '  -12.20%  '

$ws.Cells.Item(42, 5).Value = '  +4.76%  '

$c = $ws.Cells.Item(43, 4)
$c.NumberFormat = "@"
$c.Value = '393.10'
$c.Style = "Normal"
$ws.Cells.Item(43, 5).Value = '  -10.42%  '

$c = $ws.Cells.Item(44, 4)
$c.NumberFormat = "@"
$c.Value = '0.0360'
$c.Style = "Normal"
$ws.Cells.Item(44, 5).Value = '  -1.15%  '

$c = $ws.Cells.Item(45, 4)
$c.NumberFormat = "@"
$c.Value = '0.271'
$c.Style = "Normal"
$ws.Cells.Item(45, 5).Value = '  -7.11%  '

$c = $ws.Cells.Item(46, 4)
$c.NumberFormat = "@"
$c.Value = '2.734.76'
$c.Style = "Normal"
$ws.Cells.Item(46, 5).Value = '  -3.83%  '

$c = $ws.Cells.Item(47, 4)
$c.NumberFormat = "@"
$c.Value = '37.32'
$c.Style = "Normal"
$ws.Cells.Item(47, 5).Value = '  -4.17%  '

$c = $ws.Cells.Item(48, 4)
$c.NumberFormat = "@"
$c.Value = '129.56'
$c.Style = "Normal"
$ws.Cells.Item(48, 5).Value = '  -0.26%  '

$ws.Cells.Item(49, 5).Value = '  +0.04%  '

$ws.Cells.Item(50, 2).Value = 'ThetaToken'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$c = $ws.Cells.Item(50, 4)
$c.NumberFormat = "@"
$c.Value = '2.21'
$c.Style = "Normal"
$ws.Cells.Item(50, 5).Value = '  -1.06%  '

$ws.Cells.Item(51, 2).Value = 'Stellar'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$c = $ws.Cells.Item(51, 4)
$c.NumberFormat = "@"
$c.Value = '0.109'
$c.Style = "Normal"
$ws.Cells.Item(51, 5).Value = '  -0.59%  '
